$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '24.816.46'
$ws.Range('E2').Value = '  +0.78%  '
Set-TextValue 'D3' '1.705.99'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  -0.17%  '
Set-TextValue 'D5' '315.41'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  -0.19%  '
Set-TextValue 'D7' '0.4008'
$ws.Range('E7').Value = '  +2.81%  '
Set-TextValue 'D8' '0.4041'
$ws.Range('E8').Value = '  +0.32%  '
Set-TextValue 'D9' '1.476'
$ws.Range('E9').Value = '  -1.46%  '
Set-TextValue 'D10' '1.003'
$ws.Range('E10').Value = '  -0.08%  '
Set-TextValue 'D11' '53.66'
$ws.Range('E11').Value = '  +1.69%  '
Set-TextValue 'D12' '0.08807'
$ws.Range('E12').Value = '  +0.83%  '
Set-TextValue 'D13' '26.29'
$ws.Range('E13').Value = '  +5.87%  '
Set-TextValue 'D14' '7.518'
$ws.Range('E14').Value = '  -1.21%  '
Set-TextValue 'D15' '8.002'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  -0.36%  '
Set-TextValue 'D17' '1.691.60'
$ws.Range('E17').Value = '  +0.50%  '
Set-TextValue 'D18' '95.42'
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('E19').Value = '  +0.97%  '
Set-TextValue 'D20' '20.93'
$ws.Range('E20').Value = '  +5.99%  '
Set-TextValue 'D21' '7.285'
$ws.Range('E21').Value = '  +0.09%  '
Set-TextValue 'D22' '1.000'
$ws.Range('E22').Value = '  -0.27%  '
Set-TextValue 'D23' '14.44'
$ws.Range('E23').Value = '  +1.65%  '
Set-TextValue 'D24' '24.818.41'
Set-TextValue 'D25' '2.355'
$ws.Range('E25').Value = '  +0.11%  '
Set-TextValue 'D26' '2.897'
$ws.Range('E26').Value = '  -3.50%  '
Set-TextValue 'D27' '6.416'
$ws.Range('E27').Value = '  +23.09%  '
Set-TextValue 'D28' '23.10'
$ws.Range('E28').Value = '  +1.81%  '
Set-TextValue 'D29' '161.30'
$ws.Range('E29').Value = '  -0.12%  '
Set-TextValue 'D30' '143.62'
$ws.Range('E30').Value = '  +5.46%  '
Set-TextValue 'D31' '8.221'
$ws.Range('E31').Value = '  -4.26%  '
Set-TextValue 'D32' '2.278'
$ws.Range('E32').Value = '  +14.65%  '
Set-TextValue 'D33' '1.889.04'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D34' '0.03205'
$ws.Range('E34').Value = '  +10.27%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D35' '0.08648'
$ws.Range('E35').Value = '  -1.28%  '
Set-TextValue 'D36' '7.277'
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('E37').Value = '  -0.44%  '
Set-TextValue 'D38' '0.2859'
$ws.Range('E38').Value = '  +5.33%  '
Set-TextValue 'D39' '0.8410'
$ws.Range('E39').Value = '  +7.85%  '
Set-TextValue 'D40' '0.09449'
$ws.Range('E40').Value = '  +3.78%  '
Set-TextValue 'D41' '10.73'
$ws.Range('E41').Value = '  -0.36%  '
Set-TextValue 'D42' '14.26'
$ws.Range('E42').Value = '  +0.59%  '
Set-TextValue 'D43' '1.481'
$ws.Range('E43').Value = '  +1.87%  '
Set-TextValue 'D44' '17.54'
$ws.Range('E44').Value = '  +5.70%  '
Set-TextValue 'D45' '2.732'
$ws.Range('E45').Value = '  +5.98%  '
Set-TextValue 'D46' '0.7415'
$ws.Range('E46').Value = '  +3.28%  '
Set-TextValue 'D47' '4.221'
$ws.Range('E47').Value = '  +0.66%  '
Set-TextValue 'D48' '1.368'
$ws.Range('E48').Value = '  +2.51%  '
$ws.Range('E49').Value = '  -0.14%  '
Set-TextValue 'D50' '140.56'
$ws.Range('E50').Value = '  +2.14%  '
Set-TextValue 'D51' '0.08379'
$ws.Range('E51').Value = '  +5.12%  '
